$d = $word.ActiveDocument

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$pkgNs = "http://schemas.microsoft.com/office/2006/xmlPackage"

$fragments = @(
    '<w:pPr><w:spacing w:after="0"/></w:pPr>',
    '<w:r><w:br w:type="page"/></w:r>',
    '<w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>LOGROS-</w:t></w:r>',
    '<w:pPr><w:spacing w:after="0"/></w:pPr>',
    '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>1</w:t></w:r>',
    '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Registrado</w:t></w:r>',
    '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">Bienvenido a Project </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Dungeon</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r>',
    '<w:pPr><w:spacing w:after="0"/></w:pPr>',
    '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>2</w:t></w:r>',
    '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Personaliza</w:t></w:r>',
    '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Cambiemos un poco el avatar.</w:t></w:r>',
    '<w:pPr><w:spacing w:after="0"/></w:pPr>',
    '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>3</w:t></w:r>',
    '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Línea</w:t></w:r><w:r><w:t xml:space="preserve"> de partida</w:t></w:r>',
    '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Comienza una nueva partida.</w:t></w:r>',
    '<w:pPr><w:spacing w:after="0"/></w:pPr>',
    '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>4</w:t></w:r>',
    '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Nuevos horizontes</w:t></w:r>',
    '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Atraviesa el portal por primera vez.</w:t></w:r>',
    '<w:pPr><w:spacing w:after="0"/></w:pPr>',
    '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>5</w:t></w:r>',
    '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Primer asesinato</w:t></w:r>',
    '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Vence un enemigo.</w:t></w:r>',
    '<w:pPr><w:spacing w:after="0"/></w:pPr>',
    '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>6</w:t></w:r>',
    '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">Cazador de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Goblins</w:t></w:r><w:proofErr w:type="spellEnd"/>',
    '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">Mata 10 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Goblins</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r>'
)

foreach ($frag in $fragments) {
    $target = $d.Paragraphs.Last.Range
    $target.Collapse(0)
    $innerDoc = "<w:p xmlns:w=`"$wNs`">" + $frag + "</w:p>"
    $bodyXml = "<w:document xmlns:w=`"$wNs`"><w:body>" + $innerDoc + "</w:body></w:document>"
    $pkgXml = "<?xml version=`"1.0`" standalone=`"yes`"?><pkg:package xmlns:pkg=`"$pkgNs`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData>" + $bodyXml + "</pkg:xmlData></pkg:part></pkg:package>"
    [void]$target.InsertXML($pkgXml)
}

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
